# Adds two new "data_transform" model columns to the results sheet:
#   - "9. ws+sin(wd)+cos(wd)-2"   inserted before the existing column J
#   - "15. ws*sin(wd)+ws*cos(wd)-2" inserted before the (new) column P
# All existing columns from J (old) onward shift right to make room, exactly
# like using Excel's "Insert Column" on the column header.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Insert a blank column in front of column J (old "9. ws+sin(wd)+cos(wd)-3"
#    column). Everything from J..O shifts right to K..P.
# ---------------------------------------------------------------------
$ws.Range("J1").EntireColumn.Insert()

# Fill in the new column J - header + 8 data values.
$ws.Range("J1").Value = "9. ws+sin(wd)+cos(wd)-2"

$ws.Range("J2").Value = 0.036900000000000002
$ws.Range("J3").Value = 0.044999999999999998
$ws.Range("J4").Value = 0.0315
$ws.Range("J5").Value = 0.023699999999999999
$ws.Range("J6").Value = 0.038300000000000001
$ws.Range("J7").Value = 0.037499999999999999
$ws.Range("J8").Value = 0.038300000000000001
$ws.Range("J9").Value = 0.033000000000000002

# Match formatting with the left neighbour column I for rows whose label had
# the bold/italic "missing data" style applied (rows 4, 5, 9).
$ws.Range("I4").Copy()
$ws.Range("J4").PasteSpecial(-4122)
$ws.Range("I5").Copy()
$ws.Range("J5").PasteSpecial(-4122)
$ws.Range("I9").Copy()
$ws.Range("J9").PasteSpecial(-4122)

# Renumber the model labels that shifted right by one slot (old "9." .. "14."
# become "10." .. "15.").
$ws.Range("K1").Value = "10. ws+sin(wd)+cos(wd)-3"
$ws.Range("L1").Value = "11. ws*sin(wd)"
$ws.Range("M1").Value = "12. ws*cos(wd)"
$ws.Range("N1").Value = "13. ws*cos(wd)-3"
$ws.Range("O1").Value = "14. ws*sin(wd)+ws*cos(wd)"

# ---------------------------------------------------------------------
# 2) Insert a blank column in front of (the now shifted) column P - this was
#    originally column O ("14. ws*sin(wd)+ws*cos(wd)-3"), now sitting at O.
#    After this insert, the old content shifts from O to Q.
# ---------------------------------------------------------------------
$ws.Range("P1").EntireColumn.Insert()

# Fill in the new column P - header + 8 data values.
$ws.Range("P1").Value = "15. ws*sin(wd)+ws*cos(wd)-2"

$ws.Range("P2").Value = 0.045999999999999999
$ws.Range("P3").Value = 0.043700000000000003
$ws.Range("P4").Value = 0.045900000000000003
$ws.Range("P5").Value = 0.043700000000000003
$ws.Range("P6").Value = 0.041599999999999998
$ws.Range("P7").Value = 0.042599999999999999
$ws.Range("P8").Value = 0.044499999999999998
$ws.Range("P9").Value = 0.046399999999999997

# Match formatting with the left neighbour column O for row 7 (already
# carried the alternate style), and explicitly apply the same style used
# for rows 8 and 9 in that block.
$ws.Range("O7").Copy()
$ws.Range("P7").PasteSpecial(-4122)
$ws.Range("O7").Copy()
$ws.Range("P8").PasteSpecial(-4122)
$ws.Range("O7").Copy()
$ws.Range("P9").PasteSpecial(-4122)

# Final renumber: the last surviving original label ("14.") moved to Q and
# becomes "16.".
$ws.Range("Q1").Value = "16. ws*sin(wd)+ws*cos(wd)-3"

# ---------------------------------------------------------------------
# Restore the selection/view the file was left at.
# ---------------------------------------------------------------------
$ws.Range("P9").Select()
